$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40: new data for 四方坪站 (shared string "四方坪站")
$ws.Range("A40").Value = 45950
$ws.Range("B40").Value = "四方坪站"
$ws.Range("C40").Value = 10820.19
$ws.Range("D40").Value = 8805.35
$ws.Range("E40").Value = 3759.82
$ws.Range("F40").Value = 427

# Row 41: new data for 高岭站 (shared string "高岭站")
$ws.Range("A41").Value = 45950
$ws.Range("B41").Value = "高岭站"
$ws.Range("C41").Value = 4517.3
$ws.Range("D41").Value = 3687.48
$ws.Range("E41").Value = 1280.75
$ws.Range("F41").Value = 157

# Match the author's final selection state
$ws.Range("I43").Select()
